# Delete row 172 ("草を食む幸運な子羊たち" post) entirely.
# This shifts all subsequent rows (173-184) up by one, matching the diff,
# and the sheet's used range shrinks from A1:C184 to A1:C183.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(172).Delete()
